# Update Betfair Back/Lay odds for 2026-01-13
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 2.66
$ws.Range("H2").Value = 2.28
$ws.Range("I2").Value = 2.96
$ws.Range("K2").Value = 5.5
$ws.Range("P2").Value = 1.75

# Row 3
$ws.Range("F3").Value = 1.85
$ws.Range("H3").Value = 4.5
$ws.Range("J3").Value = 4.1
$ws.Range("K3").Value = 4.2
$ws.Range("P3").Value = 2.48
$ws.Range("U3").Value = 2.5

# Row 5
$ws.Range("G5").Value = 2.12
$ws.Range("J5").Value = 3.55

# Row 6
$ws.Range("F6").Value = 1.4
$ws.Range("G6").Value = 1.41
$ws.Range("I6").Value = 9.199999999999999
$ws.Range("K6").Value = 5.8
$ws.Range("N6").Value = 6.6
$ws.Range("P6").Value = 2.9
$ws.Range("Q6").Value = 1.48
$ws.Range("R6").Value = 1.76
$ws.Range("S6").Value = 2.18
$ws.Range("U6").Value = 2.22
$ws.Range("Y6").Value = 40
$ws.Range("Z6").Value = 1000
$ws.Range("AA6").Value = 290
$ws.Range("AC6").Value = 14
$ws.Range("AD6").Value = 34
$ws.Range("AE6").Value = 130
$ws.Range("AF6").Value = 11
$ws.Range("AG6").Value = 11
$ws.Range("AI6").Value = 1000
$ws.Range("AK6").Value = 14
$ws.Range("AL6").Value = 30

# Row 7
$ws.Range("G7").Value = 3.75
$ws.Range("I7").Value = 2.14
$ws.Range("J7").Value = 3.85
$ws.Range("K7").Value = 3.95
$ws.Range("N7").Value = 5
$ws.Range("P7").Value = 2.34
$ws.Range("R7").Value = 1.54
$ws.Range("S7").Value = 2.7
$ws.Range("U7").Value = 2.52
$ws.Range("X7").Value = 30
$ws.Range("Z7").Value = 970
$ws.Range("AC7").Value = 9.800000000000001
$ws.Range("AE7").Value = 21
$ws.Range("AG7").Value = 17
$ws.Range("AH7").Value = 18
$ws.Range("AO7").Value = 12

# Row 8
$ws.Range("J8").Value = 4
$ws.Range("P8").Value = 1.92
$ws.Range("Q8").Value = 2.04
$ws.Range("T8").Value = 1.98
$ws.Range("AG8").Value = 10.5
